# Regenerate orders with updated distance/size labels.
# Renames within text values (Condition, Filename_Left, Filename_Right,
# Distance, Size columns): D51->D55, D64->D69, D80->D86, S30->S31.
# (S20 and S25 are left unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$rowCount = $ur.Rows.Count
$colCount = $ur.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Text

        if ($val -is [string]) {
            if ($val -match 'D51|D64|D80|S30') {
                $newVal = $val -replace 'D51','D55' -replace 'D64','D69' -replace 'D80','D86' -replace 'S30','S31'
                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
